$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# ALC row 62
$wsALC.Range("H62").Value = 2620.5
$wsALC.Range("I62").Value = 1999.3334
$wsALC.Range("J62").Value = 2993.2
$wsALC.Range("K62").Value = 1999.3334
$wsALC.Range("L62").Value = 2993.2
$wsALC.Range("M62").Value = -1375.3334
$wsALC.Range("N62").Value = -4241.2

# ALC row 65
$wsALC.Range("H65").Value = 2620.5
$wsALC.Range("I65").Value = 1999.3334
$wsALC.Range("J65").Value = 2993.2
$wsALC.Range("K65").Value = 9996.666999999999
$wsALC.Range("L65").Value = 14966
$wsALC.Range("M65").Value = -6876.666999999999
$wsALC.Range("N65").Value = -21206

# ALC row 98
$wsALC.Range("H98").Value = 4098.75
$wsALC.Range("I98").Value = 4098.75
$wsALC.Range("J98").Value = 0
$wsALC.Range("K98").Value = 4098.75
$wsALC.Range("L98").Value = 0
$wsALC.Range("M98").Value = -2600.75

# ALC row 112
$wsALC.Range("H112").Value = 4562.727
$wsALC.Range("I112").Value = 400
$wsALC.Range("J112").Value = 4979
$wsALC.Range("K112").Value = 1200
$wsALC.Range("L112").Value = 14937
$wsALC.Range("M112").Value = -92
$wsALC.Range("N112").Value = -17153

# ALC row 113
$wsALC.Range("H113").Value = 12138.917
$wsALC.Range("I113").Value = 13866.8
$wsALC.Range("J113").Value = 3499.5
$wsALC.Range("K113").Value = 13866.8
$wsALC.Range("L113").Value = 3499.5
$wsALC.Range("M113").Value = -10612.8
$wsALC.Range("N113").Value = -10007.5

# ALC row 116
$wsALC.Range("H116").Value = 12758.9
$wsALC.Range("I116").Value = 21680.8
$wsALC.Range("J116").Value = 3837
$wsALC.Range("K116").Value = 21680.8
$wsALC.Range("L116").Value = 3837
$wsALC.Range("M116").Value = -18238.8
$wsALC.Range("N116").Value = -10721

# ALC row 118
$wsALC.Range("H118").Value = 783.3333
$wsALC.Range("I118").Value = 783.3333
$wsALC.Range("J118").Value = 0
$wsALC.Range("K118").Value = 2349.9999
$wsALC.Range("L118").Value = 0
$wsALC.Range("M118").Value = -692.9998999999998

# ALC row 122
$wsALC.Range("H122").Value = 4098.75
$wsALC.Range("I122").Value = 4098.75
$wsALC.Range("J122").Value = 0
$wsALC.Range("K122").Value = 12296.25
$wsALC.Range("L122").Value = 0
$wsALC.Range("M122").Value = -9846.25

# ALC row 125
$wsALC.Range("H125").Value = 1014.4
$wsALC.Range("I125").Value = 1350
$wsALC.Range("J125").Value = 930.5
$wsALC.Range("K125").Value = 12150
$wsALC.Range("L125").Value = 8374.5
$wsALC.Range("M125").Value = -9690
$wsALC.Range("N125").Value = -13294.5

# ALC row 132
$wsALC.Range("H132").Value = 959.55817
$wsALC.Range("I132").Value = 963.3570999999999
$wsALC.Range("J132").Value = 800
$wsALC.Range("K132").Value = 2890.0713
$wsALC.Range("L132").Value = 2400
$wsALC.Range("M132").Value = -360.0712999999996

# ALC row 138
$wsALC.Range("H138").Value = 1743.5428
$wsALC.Range("I138").Value = 1699.375
$wsALC.Range("J138").Value = 1780.7368
$wsALC.Range("K138").Value = 5098.125
$wsALC.Range("L138").Value = 5342.2104
$wsALC.Range("M138").Value = 41.875

# ALC row 139
$wsALC.Range("H139").Value = 46550
$wsALC.Range("I139").Value = 0
$wsALC.Range("J139").Value = 46550
$wsALC.Range("K139").Value = 0
$wsALC.Range("L139").Value = 46550
$wsALC.Range("N139").Value = -56830

# ALC row 140
$wsALC.Range("H140").Value = 49912.11
$wsALC.Range("I140").Value = 0
$wsALC.Range("J140").Value = 49912.11
$wsALC.Range("K140").Value = 0
$wsALC.Range("L140").Value = 49912.11
$wsALC.Range("N140").Value = -60272.11

# ALC row 141
$wsALC.Range("H141").Value = 3352.611
$wsALC.Range("I141").Value = 2306.182
$wsALC.Range("J141").Value = 4997
$wsALC.Range("K141").Value = 6918.545999999999
$wsALC.Range("L141").Value = 14991
$wsALC.Range("M141").Value = -1738.545999999999

# ARM row 32
$wsARM.Range("H32").Value = 6928.074
$wsARM.Range("I32").Value = 5108.0586
$wsARM.Range("J32").Value = 10022.1
$wsARM.Range("K32").Value = 5108.0586
$wsARM.Range("L32").Value = 10022.1
$wsARM.Range("M32").Value = -4821.0586

# ARM row 45
$wsARM.Range("H45").Value = 1689.2941
$wsARM.Range("I45").Value = 841.0909
$wsARM.Range("J45").Value = 3244.3333
$wsARM.Range("K45").Value = 841.0909
$wsARM.Range("L45").Value = 3244.3333
$wsARM.Range("M45").Value = -464.0909
$wsARM.Range("N45").Value = -3998.3333

# ARM row 61
$wsARM.Range("H61").Value = 2446.3635
$wsARM.Range("I61").Value = 1393.421
$wsARM.Range("J61").Value = 9115
$wsARM.Range("K61").Value = 1393.421
$wsARM.Range("L61").Value = 9115
$wsARM.Range("M61").Value = -1181.421

# ARM row 109
$wsARM.Range("H109").Value = 67882.664
$wsARM.Range("I109").Value = 0
$wsARM.Range("J109").Value = 67882.664
$wsARM.Range("K109").Value = 0
$wsARM.Range("L109").Value = 67882.664
$wsARM.Range("N109").Value = -70656.664

# ARM row 110
$wsARM.Range("H110").Value = 1667.5333
$wsARM.Range("I110").Value = 1221.25
$wsARM.Range("J110").Value = 2177.5715
$wsARM.Range("K110").Value = 1221.25
$wsARM.Range("L110").Value = 2177.5715
$wsARM.Range("M110").Value = 823.75

# ARM row 136
$wsARM.Range("H136").Value = 2446.3635
$wsARM.Range("I136").Value = 1393.421
$wsARM.Range("J136").Value = 9115
$wsARM.Range("K136").Value = 4180.263
$wsARM.Range("L136").Value = 27345
$wsARM.Range("M136").Value = -1630.263

# BSM row 99
$wsBSM.Range("H99").Value = 0
$wsBSM.Range("I99").Value = 0
$wsBSM.Range("J99").Value = 0
$wsBSM.Range("K99").Value = 0
$wsBSM.Range("L99").Value = 0
$wsBSM.Range("M99").ClearContents()

# BSM row 105
$wsBSM.Range("H105").Value = 1918.2413
$wsBSM.Range("I105").Value = 2031.4615
$wsBSM.Range("J105").Value = 937
$wsBSM.Range("K105").Value = 2031.4615
$wsBSM.Range("L105").Value = 937
$wsBSM.Range("M105").Value = -284.4614999999999

# BSM row 128
$wsBSM.Range("H128").Value = 200
$wsBSM.Range("I128").Value = 200
$wsBSM.Range("J128").Value = 0
$wsBSM.Range("K128").Value = 600
$wsBSM.Range("L128").Value = 0
$wsBSM.Range("M128").Value = 1890

# BSM row 134
$wsBSM.Range("H134").Value = 7678.425
$wsBSM.Range("I134").Value = 8992.357
$wsBSM.Range("J134").Value = 4612.5835
$wsBSM.Range("K134").Value = 26977.071
$wsBSM.Range("L134").Value = 13837.7505
$wsBSM.Range("M134").Value = -24442.071
$wsBSM.Range("N134").Value = -18907.7505

# CRP row 31
$wsCRP.Range("H31").Value = 2971.6667
$wsCRP.Range("I31").Value = 1355.8125
$wsCRP.Range("J31").Value = 6203.375
$wsCRP.Range("K31").Value = 1355.8125
$wsCRP.Range("L31").Value = 6203.375
$wsCRP.Range("M31").Value = -1060.8125

# CRP row 34
$wsCRP.Range("H34").Value = 2971.6667
$wsCRP.Range("I34").Value = 1355.8125
$wsCRP.Range("J34").Value = 6203.375
$wsCRP.Range("K34").Value = 1355.8125
$wsCRP.Range("L34").Value = 6203.375
$wsCRP.Range("M34").Value = -1153.8125

# CRP row 86
$wsCRP.Range("H86").Value = 1807
$wsCRP.Range("I86").Value = 1807
$wsCRP.Range("J86").Value = 0
$wsCRP.Range("K86").Value = 1807
$wsCRP.Range("L86").Value = 0
$wsCRP.Range("M86").Value = -684

# CRP row 89
$wsCRP.Range("H89").Value = 1807
$wsCRP.Range("I89").Value = 1807
$wsCRP.Range("J89").Value = 0
$wsCRP.Range("K89").Value = 9035
$wsCRP.Range("L89").Value = 0
$wsCRP.Range("M89").Value = -3419

# CRP row 92
$wsCRP.Range("H92").Value = 0
$wsCRP.Range("I92").Value = 0
$wsCRP.Range("J92").Value = 0
$wsCRP.Range("K92").Value = 0
$wsCRP.Range("L92").Value = 0
$wsCRP.Range("N92").ClearContents()

# CRP row 105
$wsCRP.Range("H105").Value = 1681
$wsCRP.Range("I105").Value = 1017.2
$wsCRP.Range("J105").Value = 5000
$wsCRP.Range("K105").Value = 1017.2
$wsCRP.Range("L105").Value = 5000
$wsCRP.Range("M105").Value = 729.8

# CUL row 108
$wsCUL.Range("H108").Value = 2130.1428
$wsCUL.Range("I108").Value = 2130.1428
$wsCUL.Range("J108").Value = 0
$wsCUL.Range("K108").Value = 6390.428400000001
$wsCUL.Range("L108").Value = 0
$wsCUL.Range("M108").Value = -3510.428400000001

# CUL row 131
$wsCUL.Range("H131").Value = 786.54
$wsCUL.Range("I131").Value = 498
$wsCUL.Range("J131").Value = 789.4545000000001
$wsCUL.Range("K131").Value = 1494
$wsCUL.Range("L131").Value = 2368.3635
$wsCUL.Range("M131").Value = 3546
$wsCUL.Range("N131").Value = -12448.3635

# CUL row 132
$wsCUL.Range("H132").Value = 1587.375
$wsCUL.Range("I132").Value = 1587.375
$wsCUL.Range("J132").Value = 0
$wsCUL.Range("K132").Value = 14286.375
$wsCUL.Range("L132").Value = 0
$wsCUL.Range("M132").Value = -11756.375

# GSM row 23
$wsGSM.Range("H23").Value = 0
$wsGSM.Range("I23").Value = 0
$wsGSM.Range("J23").Value = 0
$wsGSM.Range("K23").Value = 0
$wsGSM.Range("L23").Value = 0
$wsGSM.Range("N23").ClearContents()

# GSM row 24
$wsGSM.Range("H24").Value = 2509777.8
$wsGSM.Range("I24").Value = 10000000
$wsGSM.Range("J24").Value = 13037
$wsGSM.Range("K24").Value = 10000000
$wsGSM.Range("L24").Value = 13037
$wsGSM.Range("M24").Value = -9999827
$wsGSM.Range("N24").Value = -13383

# GSM row 102
$wsGSM.Range("H102").Value = 2025
$wsGSM.Range("I102").Value = 1933.619
$wsGSM.Range("J102").Value = 2199.4546
$wsGSM.Range("K102").Value = 1933.619
$wsGSM.Range("L102").Value = 2199.4546
$wsGSM.Range("M102").Value = -311.6189999999999

# GSM row 126
$wsGSM.Range("H126").Value = 94242.37
$wsGSM.Range("I126").Value = 3868.6667
$wsGSM.Range("J126").Value = 202690.8
$wsGSM.Range("K126").Value = 11606.0001
$wsGSM.Range("L126").Value = 608072.3999999999
$wsGSM.Range("M126").Value = -9136.000100000001
$wsGSM.Range("N126").Value = -613012.3999999999

# GSM row 132
$wsGSM.Range("H132").Value = 4699.387
$wsGSM.Range("I132").Value = 3853.52
$wsGSM.Range("J132").Value = 8223.833000000001
$wsGSM.Range("K132").Value = 11560.56
$wsGSM.Range("L132").Value = 24671.499
$wsGSM.Range("M132").Value = -9030.559999999999

# LTW row 2
$wsLTW.Range("H2").Value = 430000
$wsLTW.Range("I2").Value = 500000
$wsLTW.Range("J2").Value = 10000
$wsLTW.Range("K2").Value = 500000
$wsLTW.Range("L2").Value = 10000
$wsLTW.Range("M2").Value = -499888
$wsLTW.Range("N2").Value = -10224

# LTW row 5
$wsLTW.Range("H5").Value = 13407.333
$wsLTW.Range("I5").Value = 0
$wsLTW.Range("J5").Value = 13407.333
$wsLTW.Range("K5").Value = 0
$wsLTW.Range("L5").Value = 13407.333
$wsLTW.Range("N5").Value = -13633.333

# LTW row 7
$wsLTW.Range("H7").Value = 5974.9165
$wsLTW.Range("I7").Value = 3300.3333
$wsLTW.Range("J7").Value = 8649.5
$wsLTW.Range("K7").Value = 3300.3333
$wsLTW.Range("L7").Value = 8649.5
$wsLTW.Range("M7").Value = -3188.3333
$wsLTW.Range("N7").Value = -8873.5

# LTW row 24
$wsLTW.Range("H24").Value = 0
$wsLTW.Range("I24").Value = 0
$wsLTW.Range("J24").Value = 0
$wsLTW.Range("K24").Value = 0
$wsLTW.Range("L24").Value = 0
$wsLTW.Range("N24").ClearContents()

# LTW row 43
$wsLTW.Range("H43").Value = 10807.4
$wsLTW.Range("I43").Value = 10009
$wsLTW.Range("J43").Value = 11007
$wsLTW.Range("K43").Value = 10009
$wsLTW.Range("L43").Value = 11007
$wsLTW.Range("M43").Value = -9816
$wsLTW.Range("N43").Value = -11393

# LTW row 46
$wsLTW.Range("H46").Value = 1575.1666
$wsLTW.Range("I46").Value = 0
$wsLTW.Range("J46").Value = 1575.1666
$wsLTW.Range("K46").Value = 0
$wsLTW.Range("L46").Value = 1575.1666
$wsLTW.Range("N46").Value = -1951.1666

# LTW row 126
$wsLTW.Range("H126").Value = 5974.9165
$wsLTW.Range("I126").Value = 3300.3333
$wsLTW.Range("J126").Value = 8649.5
$wsLTW.Range("K126").Value = 9900.999899999999
$wsLTW.Range("L126").Value = 25948.5
$wsLTW.Range("M126").Value = -7430.999899999999
$wsLTW.Range("N126").Value = -30888.5

# LTW row 132
$wsLTW.Range("H132").Value = 2069.818
$wsLTW.Range("I132").Value = 1474.75
$wsLTW.Range("J132").Value = 2409.8572
$wsLTW.Range("K132").Value = 4424.25
$wsLTW.Range("L132").Value = 7229.571599999999
$wsLTW.Range("M132").Value = -1894.25

# WVR row 25
$wsWVR.Range("H25").Value = 0
$wsWVR.Range("I25").Value = 0
$wsWVR.Range("J25").Value = 0
$wsWVR.Range("K25").Value = 0
$wsWVR.Range("L25").Value = 0
$wsWVR.Range("N25").ClearContents()

# WVR row 30
$wsWVR.Range("H30").Value = 0
$wsWVR.Range("I30").Value = 0
$wsWVR.Range("J30").Value = 0
$wsWVR.Range("K30").Value = 0
$wsWVR.Range("L30").Value = 0
$wsWVR.Range("N30").ClearContents()

# WVR row 37
$wsWVR.Range("H37").Value = 67029
$wsWVR.Range("I37").Value = 0
$wsWVR.Range("J37").Value = 67029
$wsWVR.Range("K37").Value = 0
$wsWVR.Range("L37").Value = 67029
$wsWVR.Range("N37").Value = -67435

# WVR row 126
$wsWVR.Range("H126").Value = 6133.087
$wsWVR.Range("I126").Value = 5535.4375
$wsWVR.Range("J126").Value = 7499.143
$wsWVR.Range("K126").Value = 16606.3125
$wsWVR.Range("L126").Value = 22497.429
$wsWVR.Range("M126").Value = -14136.3125
$wsWVR.Range("N126").Value = -27437.429

# WVR row 132
$wsWVR.Range("H132").Value = 3157.963
$wsWVR.Range("I132").Value = 2854.0908
$wsWVR.Range("J132").Value = 4495
$wsWVR.Range("K132").Value = 8562.2724
$wsWVR.Range("L132").Value = 13485
$wsWVR.Range("M132").Value = -6032.2724
$wsWVR.Range("N132").Value = -18545

# WVR row 136
$wsWVR.Range("H136").Value = 4957.933
$wsWVR.Range("I136").Value = 5796.5
$wsWVR.Range("J136").Value = 3999.5715
$wsWVR.Range("K136").Value = 17389.5
$wsWVR.Range("L136").Value = 11998.7145
$wsWVR.Range("M136").Value = -14839.5
$wsWVR.Range("N136").Value = -17098.7145
